$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cronograma")

$newUrl = "https://www.tomsplanner.es/public/calendariodesarrollo-pqrs2?"

# Preserve B3's existing cell formatting (hyperlink style) in a scratch cell so
# re-creating the hyperlink below doesn't introduce a new/duplicate style.
$ws.Range("B3").Copy()
$ws.Range("Z100").PasteSpecial(-4122)  # xlPasteFormats

# Point the "Tom's Planner" link at its new address and update the displayed text.
$ws.Hyperlinks.Delete()
$ws.Range("B3").Value = $newUrl
$ws.Hyperlinks.Add($ws.Range("B3"), $newUrl)

# Restore the original formatting that the fresh hyperlink just clobbered.
$ws.Range("Z100").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Z100").Clear()

# Rename the calendar-tool label to match the new naming.
$ws.Range("E6").Value = "cronograma"

# The schedule sheet is now the one the user is looking at.
$ws.Activate()
$ws.Range("E7").Select()
